# Scheduled-runner data refresh: update hardcoded market-price/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) across all
# class sheets. Values below are plain numbers (no formulas involved).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4186.069
$ws.Range("I40").Value = 3966
$ws.Range("J40").Value = 4211.4614
$ws.Range("K40").Value = 3966
$ws.Range("L40").Value = 4211.4614
$ws.Range("M40").Value = -3791
$ws.Range("N40").Value = -4561.4614
$ws.Range("H43").Value = 3273.5715
$ws.Range("I43").Value = 2283.1
$ws.Range("K43").Value = 2283.1
$ws.Range("M43").Value = -2214.1
$ws.Range("H74").Value = 5562.5
$ws.Range("J74").Value = 3500
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5372
$ws.Range("H77").Value = 5562.5
$ws.Range("J77").Value = 3500
$ws.Range("L77").Value = 17500
$ws.Range("N77").Value = -26860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2056.375
$ws.Range("I2").Value = 1188.9
$ws.Range("K2").Value = 1188.9
$ws.Range("M2").Value = -1075.9
$ws.Range("H61").Value = 4117.72
$ws.Range("I61").Value = 3888.889
$ws.Range("J61").Value = 4246.4375
$ws.Range("K61").Value = 3888.889
$ws.Range("L61").Value = 4246.4375
$ws.Range("M61").Value = -3676.889
$ws.Range("N61").Value = -4670.4375
$ws.Range("H74").Value = 3843.9167
$ws.Range("I74").Value = 4125.5557
$ws.Range("J74").Value = 2999
$ws.Range("K74").Value = 4125.5557
$ws.Range("L74").Value = 2999
$ws.Range("M74").Value = -3251.5557
$ws.Range("N74").Value = -4747
$ws.Range("H77").Value = 3843.9167
$ws.Range("I77").Value = 4125.5557
$ws.Range("J77").Value = 2999
$ws.Range("K77").Value = 20627.7785
$ws.Range("L77").Value = 14995
$ws.Range("M77").Value = -16259.7785
$ws.Range("N77").Value = -23731
$ws.Range("H116").Value = 2056.375
$ws.Range("I116").Value = 1188.9
$ws.Range("K116").Value = 1188.9
$ws.Range("M116").Value = 1105.1
$ws.Range("H122").Value = 3858.9539
$ws.Range("I122").Value = 3302.2327
$ws.Range("J122").Value = 4947.091
$ws.Range("K122").Value = 9906.6981
$ws.Range("L122").Value = 14841.273
$ws.Range("M122").Value = -7456.6981
$ws.Range("N122").Value = -19741.273
$ws.Range("H132").Value = 6810.75
$ws.Range("I132").Value = 6810.75
$ws.Range("K132").Value = 20432.25
$ws.Range("M132").Value = -17902.25
$ws.Range("H136").Value = 4117.72
$ws.Range("I136").Value = 3888.889
$ws.Range("J136").Value = 4246.4375
$ws.Range("K136").Value = 11666.667
$ws.Range("L136").Value = 12739.3125
$ws.Range("M136").Value = -9116.667000000001
$ws.Range("N136").Value = -17839.3125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2056.375
$ws.Range("I3").Value = 1188.9
$ws.Range("K3").Value = 1188.9
$ws.Range("M3").Value = -1074.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 4999.5
$ws.Range("I23").Value = 4999.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 4999.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -4759.5
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 4999.5
$ws.Range("I27").Value = 4999.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 4999.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -4807.5
$ws.Range("N27").ClearContents()
$ws.Range("H58").Value = 3194.08
$ws.Range("I58").Value = 2888.1667
$ws.Range("K58").Value = 2888.1667
$ws.Range("M58").Value = -2685.1667
$ws.Range("H99").Value = 2517.6667
$ws.Range("I99").Value = 2275
$ws.Range("J99").Value = 3003
$ws.Range("K99").Value = 2275
$ws.Range("L99").Value = 3003
$ws.Range("M99").Value = -777
$ws.Range("N99").Value = -5999
$ws.Range("H105").Value = 2079.8096
$ws.Range("I105").Value = 1727.0834
$ws.Range("K105").Value = 1727.0834
$ws.Range("M105").Value = 19.91660000000002
$ws.Range("H126").Value = 2517.6667
$ws.Range("I126").Value = 2275
$ws.Range("J126").Value = 3003
$ws.Range("K126").Value = 6825
$ws.Range("L126").Value = 9009
$ws.Range("M126").Value = -4355
$ws.Range("N126").Value = -13949
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 7500
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -12570
$ws.Range("H136").Value = 3194.08
$ws.Range("I136").Value = 2888.1667
$ws.Range("K136").Value = 8664.500100000001
$ws.Range("M136").Value = -6114.500100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2506
$ws.Range("I70").Value = 2506
$ws.Range("K70").Value = 7518
$ws.Range("M70").Value = -7203
$ws.Range("H73").Value = 2506
$ws.Range("I73").Value = 2506
$ws.Range("K73").Value = 7518
$ws.Range("M73").Value = -6426
$ws.Range("H75").Value = 507.5
$ws.Range("J75").Value = 507.5
$ws.Range("L75").Value = 1522.5
$ws.Range("N75").Value = -3518.5
$ws.Range("H78").Value = 507.5
$ws.Range("J78").Value = 507.5
$ws.Range("L78").Value = 4567.5
$ws.Range("N78").Value = -14551.5
$ws.Range("H98").Value = 1636
$ws.Range("I98").Value = 2010.6
$ws.Range("J98").Value = 699.5
$ws.Range("K98").Value = 6031.799999999999
$ws.Range("L98").Value = 2098.5
$ws.Range("M98").Value = -4533.799999999999
$ws.Range("N98").Value = -5094.5
$ws.Range("H103").Value = 422.1111
$ws.Range("I103").Value = 573.25
$ws.Range("K103").Value = 1719.75
$ws.Range("M103").Value = -840.75
$ws.Range("H113").Value = 1790.3914
$ws.Range("I113").Value = 914.8570999999999
$ws.Range("J113").Value = 2173.4375
$ws.Range("K113").Value = 2744.5713
$ws.Range("L113").Value = 6520.3125
$ws.Range("M113").Value = -574.5712999999996
$ws.Range("N113").Value = -10860.3125
$ws.Range("H136").Value = 1407.2
$ws.Range("I136").Value = 1337.6666
$ws.Range("J136").Value = 2033
$ws.Range("K136").Value = 4012.9998
$ws.Range("L136").Value = 6099
$ws.Range("M136").Value = 1087.0002
$ws.Range("N136").Value = -16299

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1401.125
$ws.Range("I102").Value = 1286.0741
$ws.Range("K102").Value = 1286.0741
$ws.Range("M102").Value = 335.9259
$ws.Range("H122").Value = 5981.909
$ws.Range("I122").Value = 4964.1665
$ws.Range("K122").Value = 14892.4995
$ws.Range("M122").Value = -12442.4995
$ws.Range("H132").Value = 3677
$ws.Range("I132").Value = 3752.6667
$ws.Range("K132").Value = 11258.0001
$ws.Range("M132").Value = -8728.000100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2223.625
$ws.Range("J22").Value = 2624.75
$ws.Range("L22").Value = 2624.75
$ws.Range("N22").Value = -3214.75
$ws.Range("H27").Value = 2223.625
$ws.Range("J27").Value = 2624.75
$ws.Range("L27").Value = 2624.75
$ws.Range("N27").Value = -2838.75
$ws.Range("H46").Value = 3563.0645
$ws.Range("I46").Value = 2625.5
$ws.Range("K46").Value = 2625.5
$ws.Range("M46").Value = -2437.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 29000
$ws.Range("I33").Value = 29000
$ws.Range("K33").Value = 29000
$ws.Range("M33").Value = -28750
$ws.Range("H36").Value = 29000
$ws.Range("I36").Value = 29000
$ws.Range("K36").Value = 29000
$ws.Range("M36").Value = -28750
$ws.Range("H37").Value = 24905
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 24905
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 24905
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -25311
$ws.Range("H81").Value = 3471.5652
$ws.Range("I81").Value = 2758.75
$ws.Range("K81").Value = 5517.5
$ws.Range("M81").Value = -4456.5
$ws.Range("H84").Value = 3471.5652
$ws.Range("I84").Value = 2758.75
$ws.Range("K84").Value = 27587.5
$ws.Range("M84").Value = -22283.5
$ws.Range("H107").Value = 974.25
$ws.Range("I107").Value = 948
$ws.Range("J107").Value = 983
$ws.Range("K107").Value = 2844
$ws.Range("L107").Value = 2949
$ws.Range("M107").Value = -924
$ws.Range("N107").Value = -6789
$ws.Range("H122").Value = 4643.5
$ws.Range("I122").Value = 4901.448
$ws.Range("K122").Value = 14704.344
$ws.Range("M122").Value = -12254.344
